$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 15, shifting existing row 15 (and everything
# below it) down by one. This matches the dimension change A1:A73 -> A1:A74.
$ws.Rows("15").Insert()

# Populate the newly inserted A15 with the new "Developer Manual" dropdown
# entry that now precedes the existing "Mapping variables to outputs" entry
# (which was pushed down to row 16).
$ws.Range("A15").Value = '<ul class="dropdown-menu" aria-labelledby="dropdown-articles"><li><a class="dropdown-item" href="../articles/developer_manual.html">Developer Manual</a></li>'
